# Scheduled market-data refresh for the Siren_Profits leve-crafting workbook.
# Updates the cached Universalis price columns (currentAveragePrice*) and the
# derived Leve price/profit columns (H:N) on each job sheet to the latest pull.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 19: Unbreak My Heart / Roof Tile
$ws.Range("H19").Value = 1885.3572
$ws.Range("I19").Value = 1216.6666
$ws.Range("J19").Value = 2386.875
$ws.Range("K19").Value = 1216.6666
$ws.Range("L19").Value = 2386.875
$ws.Range("M19").Value = -1041.6666
$ws.Range("N19").Value = -2736.875
# row 53: No Accounting for Waste / Enchanted Electrum Ink
$ws.Range("H53").Value = 6141.857
$ws.Range("I53").Value = 10374.25
$ws.Range("K53").Value = 10374.25
$ws.Range("M53").Value = -9737.25
# row 76: Warding Off Temptation / Enchanted Hardsilver Ink
$ws.Range("H76").Value = 4784.3335
$ws.Range("I76").Value = 4485.2
$ws.Range("K76").Value = 4485.2
$ws.Range("M76").Value = -4170.2
# row 79: The Garden of Arcane Delights (L) / Enchanted Hardsilver Ink
$ws.Range("H79").Value = 4784.3335
$ws.Range("I79").Value = 4485.2
$ws.Range("K79").Value = 4485.2
$ws.Range("M79").Value = -3393.2
# row 80: Cleansing the Wicked Humours / Hallowed Water
$ws.Range("H80").Value = 89991.56
$ws.Range("J80").Value = 8894.833000000001
$ws.Range("L80").Value = 26684.499
$ws.Range("N80").Value = -28680.499
# row 83: Washing Away the Sins (L) / Hallowed Water
$ws.Range("H83").Value = 89991.56
$ws.Range("J83").Value = 8894.833000000001
$ws.Range("L83").Value = 80053.497
$ws.Range("N83").Value = -90037.497
# row 86: Filling in the Blanks / Enchanted Aurum Regis Ink
$ws.Range("H86").Value = 21160.541
$ws.Range("I86").Value = 2415.0588
$ws.Range("J86").Value = 66685.28999999999
$ws.Range("K86").Value = 2415.0588
$ws.Range("L86").Value = 66685.28999999999
$ws.Range("M86").Value = -1292.0588
$ws.Range("N86").Value = -68931.28999999999
# row 89: Ink into Antiquity (L) / Enchanted Aurum Regis Ink
$ws.Range("H89").Value = 21160.541
$ws.Range("I89").Value = 2415.0588
$ws.Range("J89").Value = 66685.28999999999
$ws.Range("K89").Value = 12075.294
$ws.Range("L89").Value = 333426.45
$ws.Range("M89").Value = -6459.293999999998
$ws.Range("N89").Value = -344658.45
# row 135: For Tired Minds / Grade 1 Gemsap of Intelligence
$ws.Range("H135").Value = 3603
$ws.Range("I135").Value = 3603
$ws.Range("K135").Value = 32427
$ws.Range("M135").Value = -29892
# row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 3557.8096
$ws.Range("J138").Value = 4332.32
$ws.Range("L138").Value = 12996.96
$ws.Range("N138").Value = -23276.96
# row 140: Tome for Tradition / Book of Ra'Kaznar
$ws.Range("H140").Value = 78000
$ws.Range("I140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("M140").ClearContents()

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 2: Ain't Got No Ingots / Bronze Ingot
$ws.Range("H2").Value = 4499.5
$ws.Range("I2").Value = 3999
$ws.Range("J2").Value = 5000
$ws.Range("K2").Value = 3999
$ws.Range("L2").Value = 5000
$ws.Range("M2").Value = -3886
$ws.Range("N2").Value = -5226
# row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 3514.8708
$ws.Range("I32").Value = 3695.0715
$ws.Range("K32").Value = 3695.0715
$ws.Range("M32").Value = -3408.0715
# row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 5452.3335
$ws.Range("I61").Value = 5508
$ws.Range("K61").Value = 5508
$ws.Range("M61").Value = -5296
# row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 4151.5347
$ws.Range("I74").Value = 2914.743
$ws.Range("K74").Value = 2914.743
$ws.Range("M74").Value = -2040.743
# row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 4151.5347
$ws.Range("I77").Value = 2914.743
$ws.Range("K77").Value = 14573.715
$ws.Range("M77").Value = -10205.715
# row 92: Mail It In / High Steel Scale Mail of Fending
$ws.Range("H92").Value = 157170270
$ws.Range("J92").Value = 157170270
$ws.Range("L92").Value = 157170270
$ws.Range("N92").Value = -157175262
# row 96: The Gauntlet Is Cast / High Steel Gauntlets of Fending
$ws.Range("H96").Value = 18347.334
$ws.Range("J96").Value = 18347.334
$ws.Range("L96").Value = 18347.334
$ws.Range("N96").Value = -23839.334
# row 116: No Scope / Titanbronze Ingot
$ws.Range("H116").Value = 4499.5
$ws.Range("I116").Value = 3999
$ws.Range("J116").Value = 5000
$ws.Range("K116").Value = 3999
$ws.Range("L116").Value = 5000
$ws.Range("M116").Value = -1705
$ws.Range("N116").Value = -9588
# row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 1674.9111
$ws.Range("I132").Value = 863.7105
$ws.Range("J132").Value = 6078.5713
$ws.Range("K132").Value = 2591.1315
$ws.Range("L132").Value = 18235.7139
$ws.Range("M132").Value = -61.13149999999996
$ws.Range("N132").Value = -23295.7139
# row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 5452.3335
$ws.Range("I136").Value = 5508
$ws.Range("K136").Value = 16524
$ws.Range("M136").Value = -13974

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 3: Hells Bells / Bronze Ingot
$ws.Range("H3").Value = 4499.5
$ws.Range("I3").Value = 3999
$ws.Range("J3").Value = 5000
$ws.Range("K3").Value = 3999
$ws.Range("L3").Value = 5000
$ws.Range("M3").Value = -3885
$ws.Range("N3").Value = -5228
# row 6: The Unkindest Cut / Bronze Saw
$ws.Range("H6").Value = 46700
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()
# row 13: As Above, Below / Bronze Pickaxe
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
# row 22: Riveting Run / Iron Rivets
$ws.Range("H22").Value = 334
$ws.Range("I22").Value = 281
$ws.Range("K22").Value = 281
$ws.Range("M22").Value = -108
# row 50: A Weighty Question / Mythril Sledgehammer
$ws.Range("H50").Value = 55186.332
$ws.Range("J50").Value = 55186.332
$ws.Range("L50").Value = 55186.332
$ws.Range("N50").Value = -56334.332
# row 86: Through Thick and Thin / Adamantite Nugget
$ws.Range("H86").Value = 3120.8386
$ws.Range("I86").Value = 3571.875
$ws.Range("J86").Value = 1574.4286
$ws.Range("K86").Value = 3571.875
$ws.Range("L86").Value = 1574.4286
$ws.Range("M86").Value = -2448.875
$ws.Range("N86").Value = -3820.4286
# row 89: Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws.Range("H89").Value = 3120.8386
$ws.Range("I89").Value = 3571.875
$ws.Range("J89").Value = 1574.4286
$ws.Range("K89").Value = 17859.375
$ws.Range("L89").Value = 7872.143
$ws.Range("M89").Value = -12243.375
$ws.Range("N89").Value = -19104.143

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 74: License to Heal / Dark Chestnut Rod
$ws.Range("H74").Value = 74165
$ws.Range("J74").Value = 74165
$ws.Range("L74").Value = 74165
$ws.Range("N74").Value = -75913
# row 77: Purified Polyrhythm (L) / Dark Chestnut Rod
$ws.Range("H77").Value = 74165
$ws.Range("J77").Value = 74165
$ws.Range("L77").Value = 222495
$ws.Range("N77").Value = -231231
# row 95: Standing on Ceremony / High Steel Fork
$ws.Range("H95").Value = 124826260
$ws.Range("J95").Value = 124826260
$ws.Range("L95").Value = 124826260
$ws.Range("N95").Value = -124831752
# row 107: Built to Last / White Oak Lumber
$ws.Range("H107").Value = 23133.6
$ws.Range("J107").Value = 1885.2
$ws.Range("L107").Value = 1885.2
$ws.Range("N107").Value = -5725.2
# row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 6341.3
$ws.Range("I134").Value = 3392.6365
$ws.Range("J134").Value = 9945.223
$ws.Range("K134").Value = 10177.9095
$ws.Range("L134").Value = 29835.669
$ws.Range("M134").Value = -7642.9095
$ws.Range("N134").Value = -34905.669
# row 137: Lament of the Lazylump / Dark Mahogany Fishing Rod
$ws.Range("H137").Value = 119552
$ws.Range("J137").Value = 119552
$ws.Range("L137").Value = 119552
$ws.Range("N137").Value = -129752

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 120: A Happy End / Paella
$ws.Range("H120").Value = 7404.8
$ws.Range("I120").Value = 7404.8
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 22214.4
$ws.Range("L120").Value = 0
$ws.Range("M120").Value = -17376.4
$ws.Range("N120").ClearContents()

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 80: Needs More Prayerbell / Hardsilver Ingot
$ws.Range("H80").Value = 10000.571
$ws.Range("I80").Value = 15001
$ws.Range("J80").Value = 3333.3333
$ws.Range("K80").Value = 15001
$ws.Range("L80").Value = 3333.3333
$ws.Range("M80").Value = -14003
$ws.Range("N80").Value = -5329.3333
# row 82: Appeasing the Astromancer / Hardsilver Planisphere
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
# row 83: With a Noise That Reaches Heaven (L) / Hardsilver Ingot
$ws.Range("H83").Value = 10000.571
$ws.Range("I83").Value = 15001
$ws.Range("J83").Value = 3333.3333
$ws.Range("K83").Value = 75005
$ws.Range("L83").Value = 16666.6665
$ws.Range("M83").Value = -70013
$ws.Range("N83").Value = -26650.6665
# row 85: Silver Bar of Upcycling (L) / Hardsilver Planisphere
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
# row 120: A Beneficent Elegy / Petalite Choker of Healing
$ws.Range("H120").Value = 19997
$ws.Range("J120").Value = 19997
$ws.Range("L120").Value = 19997
$ws.Range("N120").Value = -29673
# row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 3260.1875
$ws.Range("I132").Value = 2915.2856
$ws.Range("J132").Value = 5674.5
$ws.Range("K132").Value = 8745.856800000001
$ws.Range("L132").Value = 17023.5
$ws.Range("M132").Value = -6215.856800000001
$ws.Range("N132").Value = -22083.5
# row 141: Mask Maker / Black Star Mask of Casting
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 93: Hide to Go Seek / Gagana Leather
$ws.Range("H93").Value = 4927
$ws.Range("J93").Value = 985.4286
$ws.Range("L93").Value = 985.4286
$ws.Range("N93").Value = -3481.4286
# row 101: A Stitch in Time / Marid Leather Gloves of Healing
$ws.Range("H101").Value = 35272.5
$ws.Range("J101").Value = 35272.5
$ws.Range("L101").Value = 35272.5
$ws.Range("N101").Value = -41762.5
# row 136: Respect for Br'aax / Br'aax Leather
$ws.Range("H136").Value = 10598.385
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 10598.385
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 31795.155
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -36895.155

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 12: This Is Why You Can't Have Nice Things / Amateur's Breeches
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
# row 101: Who War It Better / Serge Hose of Aiming
$ws.Range("H101").Value = 15200.667
$ws.Range("J101").Value = 15301
$ws.Range("L101").Value = 15301
$ws.Range("N101").Value = -21791
# row 103: To the Tops / Serge Gambison of Healing
$ws.Range("H103").Value = 60199.5
$ws.Range("J103").Value = 60199.5
$ws.Range("L103").Value = 60199.5
$ws.Range("N103").Value = -62543.5
# row 105: One Winged Angle / Twinsilk Coat of Casting
$ws.Range("H105").Value = 40807.5
$ws.Range("J105").Value = 40807.5
$ws.Range("L105").Value = 40807.5
$ws.Range("N105").Value = -47795.5
# row 125: Color Coated / Almasty Serge Coat of Healing
$ws.Range("H125").Value = 50000
$ws.Range("J125").Value = 50000
$ws.Range("L125").Value = 50000
$ws.Range("N125").Value = -59840
# row 140: Glamorous Gloves / Thunderyards Silk Gloves of Casting
$ws.Range("H140").Value = 89982.25
$ws.Range("J140").Value = 89982.25
$ws.Range("L140").Value = 89982.25
$ws.Range("N140").Value = -100342.25

